$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 16 (shifts all existing data from row 16
# downward by 3 rows, including the last row of data).
$ws.Range("A16:A18").EntireRow.Insert()

# Populate the 3 newly inserted rows with data for days 15-17 of 08/2025.
$newData = @(
    @(15, 42971.8, 8, 2025, "08/2025"),
    @(16, 20952.5, 8, 2025, "08/2025"),
    @(17, 3877,    8, 2025, "08/2025")
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = 16 + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
